$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need the formula/T()
# + copy/paste-special-values trick below so Excel keeps them as TEXT
# (matching the source inlineStr cells) instead of silently parsing
# them into a Number and dropping significant trailing zeros.
function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $escaped = $text.Replace("""", """""")
    $cell.Value = "=T(""" + $escaped + """)"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Range("D2").Value = '64.401.64'
$ws.Range("E2").Value = '  +0.10%  '

$ws.Range("D3").Value = '3.142.68'
$ws.Range("E3").Value = '  -0.58%  '

$ws.Range("E4").Value = '  -0.01%  '

Set-TextValue "D5" '610.77'
$ws.Range("E5").Value = '  +0.68%  '

Set-TextValue "D6" '143.95'
$ws.Range("E6").Value = '  -2.65%  '

$ws.Range("D8").Value = '3.139.84'
$ws.Range("E8").Value = '  -0.53%  '

$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("E10").Value = '  -0.24%  '

Set-TextValue "D11" '5.44'
$ws.Range("E11").Value = '  -1.85%  '

$ws.Range("E12").Value = '  -0.06%  '

Set-TextValue "D13" '0.0000255'
$ws.Range("E13").Value = '  +1.68%  '

Set-TextValue "D14" '35.52'
$ws.Range("E14").Value = '  -0.67%  '

$ws.Range("D15").Value = '3.658.12'
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D17").Value = '64.342.43'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("D18").Value = '3.165.47'
$ws.Range("E18").Value = '  +0.30%  '

Set-TextValue "D19" '6.88'
$ws.Range("E19").Value = '  -1.12%  '

Set-TextValue "D20" '477.87'
$ws.Range("E20").Value = '  -0.97%  '

Set-TextValue "D21" '14.76'
$ws.Range("E21").Value = '  -0.46%  '

$ws.Range("E22").Value = '  +2.19%  '

Set-TextValue "D23" '7.84'
$ws.Range("E23").Value = '  +1.29%  '

Set-TextValue "D24" '13.70'
$ws.Range("E24").Value = '  -0.38%  '

Set-TextValue "D25" '85.39'
$ws.Range("E25").Value = '  +1.92%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").Value = '  -2.98%  '

Set-TextValue "D28" '8.57'
$ws.Range("E28").Value = '  +0.84%  '

$ws.Range("E30").Value = '  -4.72%  '

$ws.Range("E31").Value = '  +2.82%  '

$ws.Range("E32").Value = '  -0.08%  '

Set-TextValue "D33" '26.74'
$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("E34").Value = '  -3.75%  '

$ws.Range("E35").Value = '  +0.81%  '

$ws.Range("E36").Value = '  -0.68%  '

Set-TextValue "D37" '52.70'
$ws.Range("E37").Value = '  -2.88%  '

$ws.Range("E38").Value = '  +3.56%  '

Set-TextValue "D39" '456.30'
$ws.Range("E39").Value = '  +0.07%  '

Set-TextValue "D40" '3.03'
$ws.Range("E40").Value = '  +4.31%  '

Set-TextValue "D41" '0.0397'
$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("E42").Value = '  +0.63%  '

$ws.Range("E43").Value = '  -1.19%  '

$ws.Range("D44").Value = '2.881.27'
$ws.Range("E44").Value = '  +1.05%  '

$ws.Range("E45").Value = '  -1.40%  '

$ws.Range("E46").Value = '  -0.78%  '

$ws.Range("E47").Value = '  +6.09%  '

Set-TextValue "D48" '26.61'
$ws.Range("E48").Value = '  +0.65%  '

$ws.Range("E50").Value = '  +0.11%  '

Set-TextValue "D51" '121.09'
$ws.Range("E51").Value = '  +1.87%  '

